$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after row 73 (before the old row 74),
# shifting the existing rows 74-88 down to 76-90.
$ws.Rows("74:75").Insert()

# New row 74 (Americana (o) - Provincia de Limari)
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44642
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 100112021
$ws.Range("G74").Value = "Ají"
$ws.Range("H74").Value = "Americana (o)"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 25
$ws.Range("K74").Value = 21000
$ws.Range("L74").Value = 22000
$ws.Range("M74").Value = 21600
$ws.Range("N74").Value = "$/caja 25 kilos"
$ws.Range("O74").Value = "Provincia de Limarí"
$ws.Range("P74").Value = 864
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = "Hortaliza"

# New row 75 (Inferno - Región de Arica y Parinacota)
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 44642
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112021
$ws.Range("G75").Value = "Ají"
$ws.Range("H75").Value = "Inferno"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 22
$ws.Range("K75").Value = 18000
$ws.Range("L75").Value = 19000
$ws.Range("M75").Value = 18455
$ws.Range("N75").Value = "$/caja 12 kilos"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 1538
$ws.Range("Q75").Value = 12
$ws.Range("R75").Value = "Hortaliza"
